$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 542, shifting existing rows 542:612 down to 543:613
$ws.Rows(542).Insert()

# Populate the new row 542 with its data
$ws.Cells.Item(542, 1).Value2 = 9
$ws.Cells.Item(542, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(542, 3).Value2 = "Metropolitana"
$ws.Cells.Item(542, 4).Value2 = 45077
$ws.Cells.Item(542, 5).Value2 = 13
$ws.Cells.Item(542, 6).Value2 = 100112012
$ws.Cells.Item(542, 7).Value2 = "Espinaca"
$ws.Cells.Item(542, 8).Value2 = "Sin especificar"
$ws.Cells.Item(542, 9).Value2 = "Primera"
$ws.Cells.Item(542, 10).Value2 = 160
$ws.Cells.Item(542, 11).Value2 = 6000
$ws.Cells.Item(542, 12).Value2 = 7000
$ws.Cells.Item(542, 13).Value2 = 6500
$ws.Cells.Item(542, 14).Value2 = "$/cuna 10 kilos"
$ws.Cells.Item(542, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(542, 16).Value2 = 650
$ws.Cells.Item(542, 17).Value2 = 10
$ws.Cells.Item(542, 18).Value2 = "Hortaliza"
